$d = $word.ActiveDocument

# Update the Final Exam date/time: "Final Exam: 12/16, 2 pm " -> "Final Exam: 12/18 9am"
$d.Content.Find.Execute("Final Exam: 12/16, 2 pm", $true, $false, $false, $false, $false, $true, 1, $false, "Final Exam: 12/18 9am", 2)
